$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new rows (top to bottom) so the remaining rows land on the
#    correct final row numbers.
# ---------------------------------------------------------------------------

# Two new rows above the old "Main Menu" row (old row 5 -> new row 7)
$ws.Rows(5).Insert()
$ws.Rows(5).Insert()

# One new row between "Player" (now row 8) and "Visual Effects" (now row 9)
$ws.Rows(9).Insert()

# Two new rows at the very end, after "Music" (now row 14)
$ws.Rows(15).Insert()
$ws.Rows(16).Insert()

# Row insert copies the formatting of the row above into every column of the
# newly inserted blank row. The target layout only has B/C cells on rows
# 5, 6, 9, 15 and 16, so clear the stray A/D cells that come along for the
# ride (otherwise empty-but-styled cells would be serialized where the
# target has none at all).
$ws.Range("A5").Clear()
$ws.Range("D5").Clear()
$ws.Range("A6").Clear()
$ws.Range("D6").Clear()
$ws.Range("A9").Clear()
$ws.Range("D9").Clear()
$ws.Range("A15").Clear()
$ws.Range("D15").Clear()
$ws.Range("A16").Clear()
$ws.Range("D16").Clear()

# ---------------------------------------------------------------------------
# 2. Fix up cell values that moved / changed text.
# ---------------------------------------------------------------------------

# Row 5 - Level Up Menu (new)
$ws.Range("B5").Value = "Level Up Menu"
$ws.Range("C5").Value = "Using the following to make:`nGUI PRO Kit - Casual Game`nGUI PRO Kit - Fantasy RPG "

# Row 6 - In Game Timer (new)
$ws.Range("B6").Value = "In Game Timer"
$ws.Range("C6").Value = "Using the following to make:`nGUI PRO Kit - Casual Game`nGUI PRO Kit - Fantasy RPG "

# Row 7 - Main Menu (shifted down from old row 5), content unchanged

# Row 8 - Player (shifted down from old row 6); add category label
$ws.Range("A8").Value = "Game Play Assets"

# Row 9 - Visual Effects (shifted down from old row 7), content unchanged

# Row 10 - Platform Breaking Effects (new)
$ws.Range("B10").Value = "Platform Breaking Effects"
$ws.Range("C10").Value = "Will use Rayfire to pre break platforms and have the effect fo breaking."

# Row 11 - Environment (shifted down from old row 8), content unchanged
# Row 12 - Props (shifted down from old row 9), content unchanged

# Row 13 - Sound (shifted down from old row 10); add category label
$ws.Range("A13").Value = "Sound"

# Row 14 - Music (shifted down from old row 11); update the "how to find" text
$ws.Range("C14").Value = "The game itself uses speed runner music but realised probs doesn't fit theme as there need to be a reason it loops.`nCheck binding of Issac for how they do music and then select from pack Jay Ray gave us.`n"

# Row 15 - Other / Icons (new)
$ws.Range("A15").Value = "Other"
$ws.Range("B15").Value = "Icons"
$ws.Range("C15").Value = "This is the gameplay icon to play the game. This will be using:`nGUI PRO Kit - Casual Game`nGUI PRO Kit - Fantasy RPG "

# Row 16 - Curser (new)
$ws.Range("B16").Value = "Curser"
$ws.Range("C16").Value = "This is the mouse to play the game. This will be using:`nGUI PRO Kit - Casual Game`nGUI PRO Kit - Fantasy RPG "

# ---------------------------------------------------------------------------
# 3. Row heights (match target layout as closely as the engine allows).
# ---------------------------------------------------------------------------
$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).RowHeight = 43.2
$ws.Rows(10).RowHeight = 28.8
$ws.Rows(14).RowHeight = 86.4
$ws.Rows(15).RowHeight = 43.2
$ws.Rows(16).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 4. Re-point the hyperlinks, which moved from D5/D11 to D7/D14.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.gameuidatabase.com/uploads/The-End-Is-Nigh07052020-054640-96550.jpg", [Type]::Missing, [Type]::Missing, "https://www.gameuidatabase.com/uploads/The-End-Is-Nigh07052020-054640-96550.jpg")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://drive.google.com/drive/folders/1aV9m2S6_pb1bOSTbW7qG8IN_3aAU18Vu", [Type]::Missing, [Type]::Missing, "https://drive.google.com/drive/folders/1aV9m2S6_pb1bOSTbW7qG8IN_3aAU18Vu`n")

# ---------------------------------------------------------------------------
# 5. Column widths / sheet layout tweaks.
# ---------------------------------------------------------------------------
$ws.Columns("A:B").ColumnWidth = 15.666666666666666
$ws.Columns("C").ColumnWidth = 63.166666666666664
$ws.Columns("D").ColumnWidth = 133.66666666666666

# ---------------------------------------------------------------------------
# 6. View state (zoom / scroll position) - best effort.
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$win.Zoom = 70
